# Grid Bat Cap Growth per Unit Net Revenue.xlsx
# "Copy in EU-2024-develop branch" -- this copy of the US model was carried
# over into the EU model branch and re-calibrated for the EU:
#   - About sheet gets a note explaining the EU calibration source (BNEF).
#   - The GBCGpUNR input value is recalibrated from 2000 to 400.

$wb = $excel.ActiveWorkbook

$wsAbout     = $wb.Worksheets.Item("About")
$wsGBCGpUNR  = $wb.Worksheets.Item("GBCGpUNR")

# About!B9 -- new note cell next to the "Notes" label in A9.
$wsAbout.Range("B9").Value = "For the EU, roughly calibrated against EU scenario from BNEF."

# GBCGpUNR!B2 -- recalibrated input value for the EU model.
$wsGBCGpUNR.Range("B2").Value = 400

# Leave the About sheet's selection where the edit was made.
$wsAbout.Activate()
$wsAbout.Range("B10").Select()
